$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A13").Value = "Report Sabz_Erja  (Old form of 14_4)"
$ws.Range("B13").Value = "2021 May 19"
$ws.Range("C13").Value = "1400/02/29"

$ws.Range("C12").Select()
